$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the mark values for Ravi, Kush, Shyam
$ws.Range("B2").Value = 2
$ws.Range("B3").Value = 5
$ws.Range("B4").Value = 5

# Replace the AVERAGE formula with its computed value, keep conditional-format style (s=1)
$ws.Range("B5").Value = 27

# Change the highlight fill color used for the "greater than 50" condition from yellow to red
$ws.Range("B5").Interior.Color = 255
